# Updated cryptos list on Wed Apr 17 15:57:39 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) contain numeric-looking text such as
# "518.39", "60.435.57", "0.0000217" or "  -2.29%  " that must be preserved
# verbatim as text rather than being auto-coerced into floating point numbers
# by Excel. Temporarily force the whole data range to Text format before
# writing the values, then restore the default "Normal" style afterwards so
# the cells keep their original (unstyled) appearance.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

function Set-CellText($addr, $text) {
    $ws.Range($addr).Value = $text
}

# Row 2 - Bitcoin
Set-CellText "D2" "60.435.57"
Set-CellText "E2" "  -2.29%  "

# Row 3 - Ethereum
Set-CellText "D3" "2.952.31"
Set-CellText "E3" "  -2.15%  "

# Row 4 - TetherUSD
Set-CellText "E4" "  -0.08%  "

# Row 5 - BNB
Set-CellText "D5" "518.39"
Set-CellText "E5" "  -1.30%  "

# Row 6 - Solana
Set-CellText "D6" "129.13"
Set-CellText "E6" "  +0.97%  "

# Row 7 - USDC (unchanged)

# Row 8 - LidoStakedEther
Set-CellText "D8" "2.947.45"
Set-CellText "E8" "  -2.15%  "

# Row 9 - XRP
Set-CellText "D9" "0.479"
Set-CellText "E9" "  -1.72%  "

# Row 10 - Toncoin
Set-CellText "D10" "6.12"
Set-CellText "E10" "  +2.55%  "

# Row 11 - Dogecoin
Set-CellText "D11" "0.146"
Set-CellText "E11" "  -1.19%  "

# Row 12 - Cardano
Set-CellText "E12" "  -2.17%  "

# Row 13 - ShibaInu
Set-CellText "D13" "0.0000217"
Set-CellText "E13" "  -0.58%  "

# Row 14 - Avalanche
Set-CellText "D14" "32.75"
Set-CellText "E14" "  -1.21%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-CellText "D15" "3.432.98"
Set-CellText "E15" "  -1.25%  "

# Row 16 - TRON
Set-CellText "E16" "  +0.14%  "

# Row 17 - WrappedBTC
Set-CellText "D17" "60.436.08"
Set-CellText "E17" "  -2.38%  "

# Row 18 - WrappedEther
Set-CellText "D18" "2.949.04"
Set-CellText "E18" "  -2.44%  "

# Row 19 - Polkadot
Set-CellText "D19" "6.41"
Set-CellText "E19" "  -0.04%  "

# Row 20 - BitcoinCash
Set-CellText "D20" "453.77"
Set-CellText "E20" "  -3.39%  "

# Row 21 - Chainlink
Set-CellText "D21" "12.92"
Set-CellText "E21" "  -0.28%  "

# Row 22 - Polygon
Set-CellText "D22" "0.664"
Set-CellText "E22" "  -2.36%  "

# Row 23 - Uniswap
Set-CellText "D23" "6.73"
Set-CellText "E23" "  -2.49%  "

# Row 24 - Litecoin
Set-CellText "D24" "77.59"
Set-CellText "E24" "  +0.15%  "

# Row 25 - was Dai, now InternetComputer(DFINITY)
Set-CellText "B25" "InternetComputer(DFINITY)"
Set-CellText "C25" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-CellText "D25" "11.56"
Set-CellText "E25" "  -0.39%  "

# Row 26 - was InternetComputer(DFINITY), now Dai
Set-CellText "B26" "Dai"
Set-CellText "C26" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-CellText "D26" "1.00"
Set-CellText "E26" "  +0.21%  "

# Row 27 - PancakeSwap
Set-CellText "D27" "2.60"
Set-CellText "E27" "  -0.72%  "

# Row 28 - RenderToken
Set-CellText "D28" "7.58"
Set-CellText "E28" "  -4.26%  "

# Row 29 - FirstDigitalUSD
Set-CellText "D29" "0.997"
Set-CellText "E29" "  -0.29%  "

# Row 30 - EthereumClassic
Set-CellText "D30" "25.03"
Set-CellText "E30" "  -0.87%  "

# Row 31 - Mantle
Set-CellText "D31" "1.12"
Set-CellText "E31" "  +3.40%  "

# Row 32 - ImmutableX
Set-CellText "D32" "1.81"
Set-CellText "E32" "  +0.37%  "

# Row 33 - OKB
Set-CellText "D33" "54.63"
Set-CellText "E33" "  -1.68%  "

# Row 34 - Stacks
Set-CellText "D34" "2.24"
Set-CellText "E34" "  -3.02%  "

# Row 35 - NEARProtocol
Set-CellText "D35" "5.27"
Set-CellText "E35" "  +3.52%  "

# Row 36 - Filecoin
Set-CellText "D36" "5.69"
Set-CellText "E36" "  -0.99%  "

# Row 37 - Bittensor
Set-CellText "D37" "445.04"
Set-CellText "E37" "  -2.25%  "

# Row 38 - Maker
Set-CellText "D38" "3.138.86"
Set-CellText "E38" "  +3.39%  "

# Row 39 - Hedera
Set-CellText "D39" "0.0765"
Set-CellText "E39" "  -0.17%  "

# Row 40 - VeChain
Set-CellText "D40" "0.0376"
Set-CellText "E40" "  -1.71%  "

# Row 41 - Kaspa
Set-CellText "E41" "  +3.93%  "

# Row 42 - Cosmos
Set-CellText "D42" "7.89"
Set-CellText "E42" "  +0.60%  "

# Row 43 - dogwifhat
Set-CellText "D43" "2.39"
Set-CellText "E43" "  -3.36%  "

# Row 44 - USDe
Set-CellText "E44" "  +0.12%  "

# Row 45 - TheGraph
Set-CellText "D45" "0.241"
Set-CellText "E45" "  -0.42%  "

# Row 46 - InjectiveProtocol
Set-CellText "D46" "24.93"
Set-CellText "E46" "  +5.72%  "

# Row 47 - Monero
Set-CellText "D47" "118.97"
Set-CellText "E47" "  +3.37%  "

# Row 48 - Stellar
Set-CellText "E48" "  +0.94%  "

# Row 49 - was PEPE, now Fetch.AI
Set-CellText "B49" "Fetch.AI"
Set-CellText "C49" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-CellText "D49" "1.90"
Set-CellText "E49" "  -2.28%  "

# Row 50 - was Fetch.AI, now PEPE
# Build the "subscript 3" digit (U+2083) as its own string first so that
# PowerShell's "+" operator performs string concatenation instead of
# numeric addition (which happens when a char is combined with a numeric
# looking string operand).
$sub3 = [string][char]0x2083
Set-CellText "B50" "PEPE"
Set-CellText "C50" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-CellText "D50" ("0.0" + $sub3 + "0498")
Set-CellText "E50" "  -1.72%  "

# Row 51 - BitgetToken
Set-CellText "E51" "  +7.17%  "

# Restore the default (unstyled) look for the data range now that all the
# text values have been written, so no cell ends up with a lingering
# explicit Text number format.
$dataRange.Style = "Normal"
